$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and  d.tumor_grade In ["Intermediate Grade"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
 order By ss.study_subject_id ASC LIMIT 100 
'@

$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and  d.tumor_grade In ["Intermediate Grade"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`,
            lp.test_name as Platform
 order By samp.sample_id ASC LIMIT 100
'@

$filesQuery = @'
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and  d.tumor_grade In ["Intermediate Grade"]
RETURN  f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
    f.file_size AS `Size`,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`,
    ss.disease_subtype as `Diagnosis`
 order By f.file_name ASC LIMIT 100
'@

$statsQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Mucinous (colloid) Carcinoma"] and sf.grouped_recurrence_score IN ["0-5"]and d.tumor_size_group In [">4"] and  d.tumor_grade In ["Intermediate Grade"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

# Row 2: CasesTab
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("C2").Value2 = $statsQuery
$ws.Rows.Item(2).RowHeight = 360

# Row 3: SamplesTab
$ws.Range("B3").Value2 = $samplesQuery
$ws.Range("C3").Value2 = $statsQuery
$ws.Rows.Item(3).RowHeight = 409.6

# Row 4: FilesTab
$ws.Range("B4").Value2 = $filesQuery
$ws.Range("C4").Value2 = $statsQuery
$ws.Rows.Item(4).RowHeight = 316.8

# Update selection to C4 as in target
[void]$ws.Range("C4").Select()
